$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "'63.229.37"
$ws.Range("E2").Value = "  +0.67%  "
$ws.Range("D3").Value = "'2.672.81"
$ws.Range("E3").Value = "  +3.95%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'604.23"
$ws.Range("E5").Value = "  +4.11%  "
$ws.Range("D6").Value = "'143.33"
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'0.586"
$ws.Range("E8").Value = "  -0.77%  "
$ws.Range("D9").Value = "'2.675.11"
$ws.Range("E9").Value = "  +4.09%  "
$ws.Range("E10").Value = "  +0.19%  "
$ws.Range("D11").Value = "'5.68"
$ws.Range("E11").Value = "  +1.86%  "
$ws.Range("E13").Value = "  +2.14%  "
$ws.Range("D14").Value = "'27.32"
$ws.Range("E14").Value = "  +1.24%  "
$ws.Range("D15").Value = "'3.154.87"
$ws.Range("E15").Value = "  +4.10%  "
$ws.Range("D16").Value = "'63.130.90"
$ws.Range("E16").Value = "  +0.68%  "
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("E18").Value = "  +4.85%  "
$ws.Range("E19").Value = "  +3.31%  "
$ws.Range("D20").Value = "'338.91"
$ws.Range("E20").Value = "  -0.27%  "
$ws.Range("E21").Value = "  +1.42%  "
$ws.Range("D22").Value = "'6.88"
$ws.Range("E22").Value = "  +3.93%  "
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").Value = "'67.68"
$ws.Range("E24").Value = "  +0.90%  "
$ws.Range("E25").Value = "  +3.73%  "
$ws.Range("E26").Value = "  -1.31%  "
$ws.Range("E27").Value = "  -0.17%  "
$ws.Range("D28").Value = "'8.52"
$ws.Range("E28").Value = "  +3.62%  "
$ws.Range("E29").Value = "  +0.71%  "
$ws.Range("D30").Value = "'538.91"
$ws.Range("E30").Value = "  +18.65%  "
$ws.Range("D31").Value = "'7.87"
$ws.Range("E31").Value = "  -1.40%  "
$ws.Range("E32").Value = "  +5.37%  "
$ws.Range("E33").Value = "  +9.61%  "
$ws.Range("E34").Value = "  +1.81%  "
$ws.Range("D35").Value = "'173.40"
$ws.Range("E35").Value = "  -1.73%  "
$ws.Range("D36").Value = "'5.11"
$ws.Range("E36").Value = "  +14.91%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").Value = "'0.405"
$ws.Range("E38").Value = "  +1.50%  "
$ws.Range("E39").Value = "  +1.88%  "
$ws.Range("E40").Value = "  +7.41%  "
$ws.Range("D41").Value = "'175.45"
$ws.Range("E41").Value = "  +12.43%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("E43").Value = "  +1.50%  "
$ws.Range("D44").Value = "'22.07"
$ws.Range("E44").Value = "  +5.53%  "
$ws.Range("E45").Value = "  +5.27%  "
$ws.Range("D46").Value = "'0.634"
$ws.Range("E46").Value = "  +0.28%  "
$ws.Range("E47").Value = "  +0.13%  "
$ws.Range("E48").Value = "  +1.95%  "
$ws.Range("D49").Value = "'18.87"
$ws.Range("E49").Value = "  +5.27%  "
$ws.Range("D50").Value = "'1.72"
$ws.Range("E50").Value = "  +2.58%  "
$ws.Range("D51").Value = "'11.34"
$ws.Range("E51").Value = "  -0.91%  "
